$wb = $excel.ActiveWorkbook

# --- OFF sheet: row for "H" (row 2) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 374
$wsOff.Range("C2").Value = 259
$wsOff.Range("D2").Value = 52
$wsOff.Range("E2").Value = 23

# --- DEF sheet: row for "H" (row 2) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 483
$wsDef.Range("C2").Value = 325
$wsDef.Range("D2").Value = 133
$wsDef.Range("E2").Value = 58
$wsDef.Range("G2").Value = 6
